# Update "Datos actualizados" timestamp cell
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

$ws.Range("A1").Value = "Datos actualizados a 28 de Agosto de 2020 a las 02:10"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 6045068
$ws.Range("C4").Value = 44703
$ws.Range("D4").Value = 3345188
$ws.Range("E4").Value = 2515116
$ws.Range("G4").Value = 1111
$ws.Range("H4").Value = 184764

# Row 5 - Brasil
$ws.Range("B5").Value = 3764493
$ws.Range("C5").Value = 42489
$ws.Range("E5").Value = 698517
$ws.Range("G5").Value = 970
$ws.Range("H5").Value = 118726

# Row 14 - Argentina
$ws.Range("B14").Value = 380292
$ws.Range("C14").Value = 10104
$ws.Range("E14").Value = 97784
$ws.Range("G14").Value = 211
$ws.Range("H14").Value = 8050

# Row 23 - Alemania
$ws.Range("B23").Value = 240565
$ws.Range("C23").Value = 1565
$ws.Range("D23").Value = 215495
$ws.Range("E23").Value = 15711

# Row 36 - Panama
$ws.Range("B36").Value = 89982
$ws.Range("C36").Value = 900
$ws.Range("D36").Value = 63996
$ws.Range("E36").Value = 24038
$ws.Range("G36").Value = 16
$ws.Range("H36").Value = 1948

# Row 112 - Republica de Africa Central
$ws.Range("D112").Value = 1782
$ws.Range("E112").Value = 2855

# Row 113 - Montenegro
$ws.Range("B113").Value = 4558
$ws.Range("C113").Value = 59
$ws.Range("D113").Value = 3606
$ws.Range("E113").Value = 863
$ws.Range("G113").Value = 1
$ws.Range("H113").Value = 89

# Row 115 - Suazilandia
$ws.Range("B115").Value = 4433
$ws.Range("C115").Value = 46
$ws.Range("D115").Value = 3115
$ws.Range("E115").Value = 1229
$ws.Range("G115").Value = 1
$ws.Range("H115").Value = 89

# Row 118 - Surinam
$ws.Range("B118").Value = 3793
$ws.Range("C118").Value = 69
$ws.Range("D118").Value = 2893
$ws.Range("E118").Value = 834
$ws.Range("G118").Value = 4
$ws.Range("H118").Value = 66

# Row 149 - Uruguay
$ws.Range("B149").Value = 1551
$ws.Range("C149").Value = 8
$ws.Range("D149").Value = 1333
$ws.Range("E149").Value = 175
